$d = $word.ActiveDocument

# Newer Word builds silently drop the stray w:semiHidden flag on the
# built-in "Default Paragraph Font" character style the first time the
# style set is touched. Mirror that cleanup where the object model
# allows it (older/limited hosts may not expose the setter, so this is
# best-effort and must not abort the rest of the edit).
try {
    $defaultCharStyle = $d.Styles("Default Paragraph Font")
    $defaultCharStyle.Hidden = $false
} catch {
}

# Locate the "Labels: ..." paragraph that the new help text must follow.
$rng = $d.Content
$found = $rng.Find.Execute("Labels: should gene names be displayed on the network", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $labelsPara = $rng.Paragraphs(1)

    # Insert a brand-new empty paragraph right after it ...
    $labelsPara.Range.InsertParagraphAfter()

    # ... and fill it in with the documentation for the new
    # "Use Normalized Betweenness" checkbox.
    $newPara = $labelsPara.Next()
    $newPara.Range.Text = "Use Normalized Betweenness: betweenness centrality analysis is normalized to the ‘’traffic load’’ (e.g., number of interactions) between the life-cycle steps. If left unchecked the unnormalized betweenness centrality analysis is used."
}
